$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Fund sheet: selection moves from D10 to the whole column E (active cell E1)
# ---------------------------------------------------------------------------
$wsFund = $wb.Worksheets.Item("Fund")
$wsFund.Activate()
$wsFund.Range("E1:E1048576").Select()

# ---------------------------------------------------------------------------
# Deal sheet: selection moves from B2 to C1
# ---------------------------------------------------------------------------
$wsDeal = $wb.Worksheets.Item("Deal")
$wsDeal.Activate()
$wsDeal.Range("C1").Select()

# ---------------------------------------------------------------------------
# DealRequestTracker sheet: selection moves from D3 to B2:E2 (active cell B2)
# ---------------------------------------------------------------------------
$wsDRT = $wb.Worksheets.Item("DealRequestTracker")
$wsDRT.Activate()
$wsDRT.Range("B2:E2").Select()

# ---------------------------------------------------------------------------
# ToggleButtonCheck sheet: collapse the two "Column_Name1"/"Column_Name2"
# columns (F & G) into a single "Column_Name" column (F) whose values concat
# the two variants with a "<break>" separator (and "<column>" in place of
# the commas that used to separate the sub-column names).
# ---------------------------------------------------------------------------
$wsTBC = $wb.Worksheets.Item("ToggleButtonCheck")
$wsTBC.Activate()

# Drop column G entirely - its data is folded into column F below.
$wsTBC.Range("G1").EntireColumn.Delete()

# Header for the merged column.
$wsTBC.Range("F1").Value = "Column_Name"

# Column E (ToggleButton) values: comma -> "<break>" between the two toggle
# states. Written top-to-bottom first to match authoring order.
$wsTBC.Range("E2").Value = "Fund Investments<break>Co-Investments"
$wsTBC.Range("E3").Value = "Open Questions<break>Closed"
$wsTBC.Range("E4").Value = "Third Party Event<break>Our Events"

# Column F values: merge what used to be F (Column_Name1) and G
# (Column_Name2), commas between field names become "<column>", and the two
# merged variants are joined with "<break>".
$wsTBC.Range("F2").Value = "Legal Entity<column>Fund<column>Commitment Amount(M)<column>Commitment Date<break>Legal Entity<column>Asset<column>Commitment Amount(M)<column>Commitment Date"
$wsTBC.Range("F3").Value = "Request Tracker ID<column>Date Requested<column>Request<column>Status<break>Request Tracker ID<column>Date Requested<column>Request"
$wsTBC.Range("F4").Value = "Name<column>Title<column>Email<break>Staff Name<column>Title<column>Mobile Phone"

# Widen column F now that it holds the combined content.
$wsTBC.Columns("F:F").ColumnWidth = 126.6

# View was scrolled right (topLeftCell D1) and left on an empty cell below
# the data (F14) after the edits were made.
$win = $excel.ActiveWindow
$win.ScrollColumn = 4
$win.ScrollRow = 1
$wsTBC.Range("F14").Select()
